$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two occurrence records on row 2 and row 3 had several of their
# observation-specific fields swapped between them (Id, Antal, Ost, Nord).
# Re-create that swap.

# --- Column A (Id) - numeric ---
$a2 = $ws.Range("A2").Value2
$a3 = $ws.Range("A3").Value2
$ws.Range("A2").Value2 = $a3
$ws.Range("A3").Value2 = $a2

# --- Column I (Antal) - stored as text, keep it text after the swap ---
$i2 = $ws.Range("I2").Text
$i3 = $ws.Range("I3").Text
$ws.Range("I2").Value = "'" + $i3
$ws.Range("I3").Value = "'" + $i2

# --- Column Q (Ost) - numeric ---
$q2 = $ws.Range("Q2").Value2
$q3 = $ws.Range("Q3").Value2
$ws.Range("Q2").Value2 = $q3
$ws.Range("Q3").Value2 = $q2

# --- Column R (Nord) - numeric ---
$r2 = $ws.Range("R2").Value2
$r3 = $ws.Range("R3").Value2
$ws.Range("R2").Value2 = $r3
$ws.Range("R3").Value2 = $r2
